# Daily attendance processing - 2026-01-12 06:09:55
# Rotates the comma-separated "Recorded By" list in column G so that the
# last entry moves to the front (right-rotate by one), for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
